$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.291.64"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.816.71"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.44"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.30"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.817.37"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.03"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.459.67"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.815.62"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.254.54"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.47"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.11"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.88"
$ws.Range("E22").Value = "  -5.04%  "
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.09"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.966.07"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.49"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  -5.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.43"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.781.13"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  +10.83%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  -3.44%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -4.81%  "
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.75"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  +8.23%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "415.19"
$ws.Range("E48").Value = "  -5.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.96"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.49"
$ws.Range("E51").Value = "  -1.44%  "
